# Add a new column R (year 2021) to the worksheet, mirroring the
# formatting of column Q (year 2020), then fill in the 2021 data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column Q (rows 4-14) into column R so the
# new cells pick up the same number formats / styles as their neighbours.
# Using Copy + Insert (rather than a plain Copy/paste into the existing
# blank column) so the engine preserves the exact source cell style
# (including any redundant-but-harmless formatting flags) instead of
# re-resolving to a different, visually-identical style index.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").Insert(-4161)

# Now overwrite the copied values with the actual 2021 data.
$ws.Range("R4").Value  = 2021
$ws.Range("R5").Value  = 111.17903216128188
$ws.Range("R6").Value  = 113.69236134930286
$ws.Range("R7").Value  = 114.88854111210361
$ws.Range("R8").Value  = 110.91060220352473
$ws.Range("R9").Value  = 113.02233875668462
$ws.Range("R10").Value = 110.66816227588356
$ws.Range("R11").Value = 111.40708764208969
$ws.Range("R12").Value = 109.49389157333138
$ws.Range("R13").Value = 110.97185980126036
$ws.Range("R14").Value = 110.008558587758

# Update the selected/active cell as recorded in the saved view state.
$ws.Range("T6").Select()
